$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("R1").Range("G2").Value = "3947:35:00"
$wb.Worksheets.Item("R1").Range("G3").Value = "87:07:38"
$wb.Worksheets.Item("R1").Range("G4").Value = "110:07:38"

$wb.Worksheets.Item("R2").Range("G2").Value = "12128:58:41"
$wb.Worksheets.Item("R2").Range("G3").Value = "3258:42:10"
$wb.Worksheets.Item("R2").Range("G4").Value = "496:53:44"

$wb.Worksheets.Item("R4").Range("G2").Value = "2974:48:30"
$wb.Worksheets.Item("R4").Range("G3").Value = "202:00:45"
$wb.Worksheets.Item("R4").Range("G4").Value = "90:13:10"
$wb.Worksheets.Item("R4").Range("G5").Value = "87:50:43"

$wb.Worksheets.Item("R5").Range("G2").Value = "448:47:29"

$wb.Worksheets.Item("R6").Range("G2").Value = "89:19:47"
